# Adds a debug-notes test: "Reason" / "Solution" notes about why certain
# companies could not be matched, plus a brand-new row (O'Reilly Automotive)
# for the "could not find info on macrotrends" error case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (D1:F1) -------------------------------------------
$ws.Range("D1").Value = "My notes"
$ws.Range("E1").Value = "Reason"
$ws.Range("F1").Value = "Solution"

# --- Per-row debug notes ---------------------------------------------------
$ws.Range("D25").Value = "Not publicly traded"

$ws.Range("D26").Value = "Not publicly traded, also called (edward jones)"

$ws.Range("D27").Value = "they were acquired by wesco international 'WCC'"

$ws.Range("D28").Value = "Can be found through 'erie' on marketwatch then scanning page for ' erie ' word"
$ws.Range("E28").Value = "Erie Insurance Group is the parent of erie indemnity company"
$ws.Range("F28").Value = "Use backup algorithm on market watch or use backup website"

$ws.Range("D29").Value = "They are a mutual insurance company, so not publicly traded"

$ws.Range("D30").Value = "Just 'newmont' works on marketwatch"

$ws.Range("D31").Value = "Not publicly traded"

$ws.Range("D32").Value = "Could be found using 'Packaging Corp.' and 'Packaging' on marketwatch"

$ws.Range("D33").Value = "They were purchased by Cleveland Cliffs"
$ws.Range("E33").Value = "They were purchased by Cleveland Cliffs"

$ws.Range("D34").Value = "Could not be found with both words, but could be found with just 'Polaris' on marketwatch as Polaris Inc."
$ws.Range("E34").Value = "They changed their name to Polaris Inc. in July 2019"

# --- New row 35: O'Reilly Automotive --------------------------------------
$ws.Range("A35").Value = 329
$ws.Range("B35").Value = "O'Reilly Automotive"
$ws.Range("C35").Value = "could not find info on macrotrends"
$ws.Range("D35").Value = "Can't use macrotrends, need backup website, only work on backup if there are lots of this error"

# --- Column widths (best fit) ---------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 4.0833333333333335
$ws.Columns.Item(2).ColumnWidth = 33.25
$ws.Columns.Item(3).ColumnWidth = 31.416666666666668
$ws.Columns.Item(4).ColumnWidth = 95.08333333333333
$ws.Columns.Item(5).ColumnWidth = 55.916666666666664
$ws.Columns.Item(6).ColumnWidth = 56.083333333333336

# --- Selection --------------------------------------------------------------
$ws.Range("D24").Select()
